$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9166777729988098
$ws.Range("B1").Value = 1.557201862335205
$ws.Range("C1").Value = 3.071535348892212
$ws.Range("D1").Value = 3.874866724014282
$ws.Range("E1").Value = 0.3883232772350311
